$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $text
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" '27.247.29'
Set-TextValue "E2" '  -3.35%  '

# Row 3
Set-TextValue "D3" '1.733.75'
Set-TextValue "E3" '  -3.53%  '

# Row 4
Set-TextValue "D4" '1.011'
Set-TextValue "E4" '  +0.64%  '

# Row 5
Set-TextValue "D5" '321.95'
Set-TextValue "E5" '  -4.97%  '

# Row 6
Set-TextValue "D6" '1.006'
Set-TextValue "E6" '  +0.55%  '

# Row 7
Set-TextValue "D7" '0.4223'
Set-TextValue "E7" '  -11.56%  '

# Row 8
Set-TextValue "D8" '0.3572'
Set-TextValue "E8" '  -3.50%  '

# Row 9
Set-TextValue "D9" '44.80'
Set-TextValue "E9" '  -1.56%  '

# Row 10
Set-TextValue "D10" '1.111'
Set-TextValue "E10" '  -2.76%  '

# Row 11
Set-TextValue "D11" '0.07316'
Set-TextValue "E11" '  -4.61%  '

# Row 12
Set-TextValue "D12" '1.005'
Set-TextValue "E12" '  +0.33%  '

# Row 13
Set-TextValue "D13" '21.35'
Set-TextValue "E13" '  -5.40%  '

# Row 14
Set-TextValue "D14" '6.036'
Set-TextValue "E14" '  -4.14%  '

# Row 15
Set-TextValue "D15" '7.149'
Set-TextValue "E15" '  -2.12%  '

# Row 16
Set-TextValue "D16" '1.744.69'
Set-TextValue "E16" '  -2.93%  '

# Row 17
Set-TextValue "D17" '0.00001051'
Set-TextValue "E17" '  -3.95%  '

# Row 18
Set-TextValue "D18" '83.94'
Set-TextValue "E18" '  +2.35%  '

# Row 19
Set-TextValue "D19" '0.05946'
Set-TextValue "E19" '  -11.47%  '

# Row 20
Set-TextValue "D20" '1.007'
Set-TextValue "E20" '  +0.63%  '

# Row 21
Set-TextValue "D21" '16.64'
Set-TextValue "E21" '  -4.15%  '

# Row 22
Set-TextValue "D22" '6.007'
Set-TextValue "E22" '  -6.19%  '

# Row 23
Set-TextValue "D23" '27.356.66'
Set-TextValue "E23" '  -2.98%  '

# Row 24
Set-TextValue "D24" '11.19'
Set-TextValue "E24" '  -6.76%  '

# Row 25
Set-TextValue "D25" '2.406'
Set-TextValue "E25" '  -0.01%  '

# Row 26
Set-TextValue "D26" '19.79'
Set-TextValue "E26" '  -4.23%  '

# Row 27
Set-TextValue "D27" '149.22'
Set-TextValue "E27" '  -0.96%  '

# Row 28
Set-TextValue "D28" '2.318'
Set-TextValue "E28" '  -3.55%  '

# Row 29
Set-TextValue "D29" '1.943.23'
Set-TextValue "E29" '  -3.00%  '

# Row 30
Set-TextValue "D30" '1.273'
Set-TextValue "E30" '  +0.29%  '

# Row 31
Set-TextValue "D31" '125.54'
Set-TextValue "E31" '  -6.28%  '

# Row 32
Set-TextValue "D32" '3.734'
Set-TextValue "E32" '  -7.80%  '

# Row 33
Set-TextValue "D33" '0.08979'
Set-TextValue "E33" '  -7.03%  '

# Row 34
Set-TextValue "D34" '5.486'
Set-TextValue "E34" '  -7.46%  '

# Row 35
Set-TextValue "D35" '12.23'
Set-TextValue "E35" '  +0.69%  '

# Row 36
Set-TextValue "E36" '  -2.64%  '

# Row 37
Set-TextValue "B37" 'Hedera'
Set-TextValue "C37" 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D37" '0.06047'
Set-TextValue "E37" '  -4.15%  '

# Row 38
Set-TextValue "B38" 'TheSandbox'
Set-TextValue "C38" 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue "D38" '0.6397'
Set-TextValue "E38" '  -4.08%  '

# Row 39
Set-TextValue "B39" 'VeChain'
Set-TextValue "C39" 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue "D39" '0.02240'
Set-TextValue "E39" '  -5.36%  '

# Row 40
Set-TextValue "D40" '4.942'
Set-TextValue "E40" '  -5.50%  '

# Row 41
Set-TextValue "D41" '1.172'
Set-TextValue "E41" '  -4.05%  '

# Row 42
Set-TextValue "D42" '1.007'
Set-TextValue "E42" '  +0.63%  '

# Row 43
Set-TextValue "D43" '1.412'
Set-TextValue "E43" '  -4.99%  '

# Row 44
Set-TextValue "D44" '7.778'
Set-TextValue "E44" '  -4.08%  '

# Row 45
Set-TextValue "D45" '13.42'
Set-TextValue "E45" '  -5.36%  '

# Row 46
Set-TextValue "D46" '3.745'
Set-TextValue "E46" '  -3.25%  '

# Row 47
Set-TextValue "E47" '  -5.24%  '

# Row 48
Set-TextValue "D48" '123.30'
Set-TextValue "E48" '  -4.34%  '

# Row 49
Set-TextValue "D49" '1.918'
Set-TextValue "E49" '  -6.13%  '

# Row 50
Set-TextValue "D50" '0.06791'
Set-TextValue "E50" '  -4.31%  '

# Row 51
Set-TextValue "D51" '1.088'
Set-TextValue "E51" '  -6.95%  '
